$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 288901.06
$ws.Range("I141").Value = 1206.375
$ws.Range("K141").Value = 3619.125
$ws.Range("M141").Value = 1560.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1325.44
$ws.Range("I32").Value = 1207.2529
$ws.Range("J32").Value = 2116.3845
$ws.Range("K32").Value = 1207.2529
$ws.Range("L32").Value = 2116.3845
$ws.Range("M32").Value = -920.2529
$ws.Range("N32").Value = -2690.3845
$ws.Range("H135").Value = 28000
$ws.Range("J135").Value = 28000
$ws.Range("L135").Value = 28000
$ws.Range("N135").Value = -38140
$ws.Range("H139").Value = 25000
$ws.Range("J139").Value = 25000
$ws.Range("L139").Value = 25000
$ws.Range("N139").Value = -35280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 30780
$ws.Range("J133").Value = 30780
$ws.Range("L133").Value = 30780
$ws.Range("N133").Value = -40900
$ws.Range("H134").Value = 2281.359
$ws.Range("I134").Value = 1799.1333
$ws.Range("K134").Value = 5397.3999
$ws.Range("M134").Value = -2862.3999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 27758
$ws.Range("J28").Value = 27758
$ws.Range("L28").Value = 27758
$ws.Range("N28").Value = -28248
$ws.Range("H31").Value = 3712.4243
$ws.Range("I31").Value = 2194.842
$ws.Range("J31").Value = 5772
$ws.Range("K31").Value = 2194.842
$ws.Range("L31").Value = 5772
$ws.Range("M31").Value = -1899.842
$ws.Range("N31").Value = -6362
$ws.Range("H34").Value = 3712.4243
$ws.Range("I34").Value = 2194.842
$ws.Range("J34").Value = 5772
$ws.Range("K34").Value = 2194.842
$ws.Range("L34").Value = 5772
$ws.Range("M34").Value = -1992.842
$ws.Range("N34").Value = -6176
$ws.Range("H58").Value = 7938608.5
$ws.Range("I58").Value = 1058.0454
$ws.Range("J58").Value = 26320304
$ws.Range("K58").Value = 1058.0454
$ws.Range("L58").Value = 26320304
$ws.Range("M58").Value = -855.0454
$ws.Range("N58").Value = -26320710
$ws.Range("H99").Value = 2463.3635
$ws.Range("I99").Value = 1457
$ws.Range("K99").Value = 1457
$ws.Range("M99").Value = 41
$ws.Range("H126").Value = 2463.3635
$ws.Range("I126").Value = 1457
$ws.Range("K126").Value = 4371
$ws.Range("M126").Value = -1901
$ws.Range("H134").Value = 1525
$ws.Range("I134").Value = 719.2222
$ws.Range("J134").Value = 3198.5386
$ws.Range("K134").Value = 2157.6666
$ws.Range("L134").Value = 9595.6158
$ws.Range("M134").Value = 377.3334
$ws.Range("N134").Value = -14665.6158
$ws.Range("H136").Value = 7938608.5
$ws.Range("I136").Value = 1058.0454
$ws.Range("J136").Value = 26320304
$ws.Range("K136").Value = 3174.1362
$ws.Range("L136").Value = 78960912
$ws.Range("M136").Value = -624.1361999999999
$ws.Range("N136").Value = -78966012
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3339.3845
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 3576
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 10728
$ws.Range("M9").Value = -1276
$ws.Range("N9").Value = -11176
$ws.Range("H16").Value = 1812.75
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 2780.4
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 8341.200000000001
$ws.Range("M16").Value = -427
$ws.Range("N16").Value = -8687.200000000001
$ws.Range("H20").Value = 3483.3333
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 3580
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10740
$ws.Range("M20").Value = -8773
$ws.Range("N20").Value = -11194
$ws.Range("H22").Value = 9450.5
$ws.Range("J22").Value = 9450.5
$ws.Range("L22").Value = 28351.5
$ws.Range("N22").Value = -28689.5
$ws.Range("H27").Value = 9450.5
$ws.Range("J27").Value = 9450.5
$ws.Range("L27").Value = 28351.5
$ws.Range("N27").Value = -28555.5
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -9566
$ws.Range("H138").Value = 4140.125
$ws.Range("J138").Value = 6795
$ws.Range("L138").Value = 20385
$ws.Range("N138").Value = -30665

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1564.9836
$ws.Range("I132").Value = 937.5789
$ws.Range("J132").Value = 2601.5652
$ws.Range("K132").Value = 2812.7367
$ws.Range("L132").Value = 7804.6956
$ws.Range("M132").Value = -282.7366999999999
$ws.Range("N132").Value = -12864.6956

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1692.0714
$ws.Range("I46").Value = 431.5
$ws.Range("J46").Value = 2637.5
$ws.Range("K46").Value = 431.5
$ws.Range("L46").Value = 2637.5
$ws.Range("M46").Value = -243.5
$ws.Range("N46").Value = -3013.5
$ws.Range("H132").Value = 3051.4075
$ws.Range("I132").Value = 2092.4614
$ws.Range("J132").Value = 3941.8572
$ws.Range("K132").Value = 6277.3842
$ws.Range("L132").Value = 11825.5716
$ws.Range("M132").Value = -3747.3842
$ws.Range("N132").Value = -16885.5716
$ws.Range("H133").Value = 29625
$ws.Range("J133").Value = 29625
$ws.Range("L133").Value = 29625
$ws.Range("N133").Value = -34685
$ws.Range("H139").Value = 43782.145
$ws.Range("J139").Value = 43782.145
$ws.Range("L139").Value = 43782.145
$ws.Range("N139").Value = -54062.145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1131.25
$ws.Range("I100").Value = 1083.8334
$ws.Range("K100").Value = 2167.6668
$ws.Range("M100").Value = -1626.6668
$ws.Range("H113").Value = 2049.1667
$ws.Range("J113").Value = 3071.1428
$ws.Range("L113").Value = 9213.428400000001
$ws.Range("N113").Value = -13553.4284
$ws.Range("H132").Value = 14048.489
$ws.Range("I132").Value = 3062
$ws.Range("J132").Value = 41092.152
$ws.Range("K132").Value = 9186
$ws.Range("L132").Value = 123276.456
$ws.Range("M132").Value = -6656
$ws.Range("N132").Value = -128336.456
